{"js": "// \"2018/12/10 tidy my note\"\n// The note's git cheat-sheet ended with:\n//   git add .\n//   git commit \u2013m\n//   git push        (carries the _GoBack bookmark)\n//\n// It becomes:\n//   git add .  //\u63d0\u4ea4\u65b0\u589e\u3001\u4fee\u6539\u6587\u4ef6\n//   git add \u2013A  //\u63d0\u4ea4\u65b0\u589e\u3001\u4fee\u6539\u3001\u5220\u9664\u6587\u4ef6   (new paragraph, now carries _GoBack)\n//   git commit \u2013m\n//   git push\n\nlet paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// 1. Locate the \"git add .\" paragraph (use the last match \u2014 the doc reuses\n//    similar snippets earlier, the one we want is right before \"git commit \u2013m\").\nlet addIndex = -1;\nfor (let i = paras.items.length - 1; i >= 0; i--) {\n  if (paras.items[i].text === \"git add .\") {\n    addIndex = i;\n    break;\n  }\n}\nif (addIndex === -1) {\n  throw new Error(\"Could not locate the 'git add .' paragraph\");\n}\nconst addPara = paras.items[addIndex];\n\n// 2. Append the explanatory comment to the end of that paragraph.\naddPara\n  .getRange(Word.RangeLocation.end)\n  .insertText(\"  //\u63d0\u4ea4\u65b0\u589e\u3001\u4fee\u6539\u6587\u4ef6\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Insert a new paragraph right after it for the \"git add \u2013A\" command.\naddPara.insertParagraph(\n  \"git add \u2013A  //\u63d0\u4ea4\u65b0\u589e\u3001\u4fee\u6539\u3001\u5220\u9664\u6587\u4ef6\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 4. Re-fetch paragraphs to get a fresh/live reference to the paragraph we\n//    just inserted (reusing the object returned by insertParagraph can give\n//    stale range positions for bookmark insertion).\nparas = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\nconst newPara = paras.items[addIndex + 1];\n\n// 5. Move the \"_GoBack\" bookmark off of the \"git push\" paragraph and onto\n//    the end of this newly inserted paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nnewPara.getRange(Word.RangeLocation.end).insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"2018/12/10 tidy my note\"\n# The note's git cheat-sheet ended with:\n#   git add .\n#   git commit \u2013m\n#   git push        (carries the _GoBack bookmark)\n#\n# It becomes:\n#   git add .  //\u63d0\u4ea4\u65b0\u589e\u3001\u4fee\u6539\u6587\u4ef6\n#   git add \u2013A  //\u63d0\u4ea4\u65b0\u589e\u3001\u4fee\u6539\u3001\u5220\u9664\u6587\u4ef6   (new paragraph, now carries _GoBack)\n#   git commit \u2013m\n#   git push\n\n$d = $word.ActiveDocument\n\n# 1. Locate the \"git add .\" paragraph. Walk from the end of the document\n#    since similar snippets appear earlier in this note; the one we want is\n#    the one immediately preceding \"git commit \u2013m\" / \"git push\".\n$addIndex = -1\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n    if ($t -eq \"git add .\") {\n        $addIndex = $i\n        break\n    }\n}\nif ($addIndex -eq -1) {\n    throw \"Could not find the 'git add .' paragraph\"\n}\n\n# 2. Append the explanatory comment to the end of that paragraph's text\n#    (i.e. just before its paragraph mark).\n$addPara = $d.Paragraphs.Item($addIndex)\n$endPoint = $d.Range($addPara.Range.End - 1, $addPara.Range.End - 1)\n$endPoint.InsertAfter(\"  //\u63d0\u4ea4\u65b0\u589e\u3001\u4fee\u6539\u6587\u4ef6\")\n\n# 3. Insert a brand new (empty) paragraph right after it.\n$addPara = $d.Paragraphs.Item($addIndex)\n$addContent = $d.Range($addPara.Range.Start, $addPara.Range.End - 1)\n$addContent.InsertParagraphAfter()\n\n# 4. Fill in the text of the newly created paragraph.\n$newIndex = $addIndex + 1\n$newPara = $d.Paragraphs.Item($newIndex)\n$newContent = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)\n$newContent.Text = \"git add \u2013A  //\u63d0\u4ea4\u65b0\u589e\u3001\u4fee\u6539\u3001\u5220\u9664\u6587\u4ef6\"\n\n# 5. Move the \"_GoBack\" bookmark off of the \"git push\" paragraph and onto the\n#    end of this new paragraph.\n#\n#    Note: this host (like Word itself) will not directly place a bookmark at\n#    a collapsed range that sits exactly \"at end of paragraph content, right\n#    before the pilcrow\" via Bookmarks.Add \u2014 it snaps to a bogus position.\n#    The reliable fix is the classic trick: type a one-character placeholder\n#    after the text, drop the bookmark immediately before that placeholder\n#    (a perfectly ordinary, non-edge position), then delete the placeholder.\n#    Being a stable point, the bookmark stays exactly where we put it.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$newPara = $d.Paragraphs.Item($newIndex)\n$placeholderPoint = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)\n$placeholderPoint.InsertAfter(\"X\")\n\n$newPara = $d.Paragraphs.Item($newIndex)\n$bookmarkPos = $newPara.Range.End - 2\n$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n$newPara = $d.Paragraphs.Item($newIndex)\n$placeholderRange = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)\n$placeholderRange.Delete()\n"}
